$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a range to store a plain-text value, even if the text
# looks like a number (e.g. "1.001"), while keeping the default "Normal"
# style (no lingering text-format / quote-prefix style index).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '24.024.09'
$ws.Range('E2').Value = '  -0.42%  '

$ws.Range('D3').Value = '1.663.29'
$ws.Range('E3').Value = '  +1.50%  '

Set-TextValue $ws.Range('D4') '1.001'
$ws.Range('E4').Value = '  -0.06%  '

Set-TextValue $ws.Range('D5') '309.97'
$ws.Range('E5').Value = '  +0.27%  '

Set-TextValue $ws.Range('D6') '1.001'
$ws.Range('E6').Value = '  -0.28%  '

Set-TextValue $ws.Range('D7') '0.3908'
$ws.Range('E7').Value = '  -0.80%  '

Set-TextValue $ws.Range('D8') '0.3882'
$ws.Range('E8').Value = '  +0.31%  '

Set-TextValue $ws.Range('D9') '51.67'
$ws.Range('E9').Value = '  +2.85%  '

Set-TextValue $ws.Range('D10') '1.375'
$ws.Range('E10').Value = '  +0.61%  '

$ws.Range('E11').Value = '  -0.01%  '

Set-TextValue $ws.Range('D12') '0.08503'
$ws.Range('E12').Value = '  -0.70%  '

Set-TextValue $ws.Range('D13') '24.20'
$ws.Range('E13').Value = '  +2.17%  '

Set-TextValue $ws.Range('D14') '7.253'
$ws.Range('E14').Value = '  +2.43%  '

Set-TextValue $ws.Range('D15') '8.004'
$ws.Range('E15').Value = '  +6.86%  '

Set-TextValue $ws.Range('D16') '0.00001320'
$ws.Range('E16').Value = '  +2.61%  '

$ws.Range('D17').Value = '1.659.19'
$ws.Range('E17').Value = '  +1.16%  '

Set-TextValue $ws.Range('D18') '95.13'
$ws.Range('E18').Value = '  +1.28%  '

Set-TextValue $ws.Range('D19') '0.06991'
$ws.Range('E19').Value = '  +0.97%  '

Set-TextValue $ws.Range('D20') '20.00'
$ws.Range('E20').Value = '  -1.82%  '

Set-TextValue $ws.Range('D21') '7.018'
$ws.Range('E21').Value = '  +1.46%  '

Set-TextValue $ws.Range('D22') '1.000'
$ws.Range('E22').Value = '  -0.28%  '

Set-TextValue $ws.Range('D23') '13.75'
$ws.Range('E23').Value = '  +1.10%  '

$ws.Range('D24').Value = '24.023.78'
$ws.Range('E24').Value = '  -0.38%  '

Set-TextValue $ws.Range('D25') '3.164'
$ws.Range('E25').Value = '  +9.51%  '

Set-TextValue $ws.Range('D26') '2.489'
$ws.Range('E26').Value = '  +3.45%  '

Set-TextValue $ws.Range('D27') '22.31'
$ws.Range('E27').Value = '  +0.26%  '

Set-TextValue $ws.Range('D28') '154.50'
$ws.Range('E28').Value = '  -2.06%  '

Set-TextValue $ws.Range('D29') '140.25'
$ws.Range('E29').Value = '  +0.22%  '

Set-TextValue $ws.Range('D30') '5.299'
$ws.Range('E30').Value = '  +0.59%  '

Set-TextValue $ws.Range('D31') '7.851'
$ws.Range('E31').Value = '  -3.26%  '

Set-TextValue $ws.Range('D32') '2.482'
$ws.Range('E32').Value = '  -0.25%  '

$ws.Range('D33').Value = '1.842.25'
$ws.Range('E33').Value = '  +1.01%  '

Set-TextValue $ws.Range('D34') '1.050'
$ws.Range('E34').Value = '  +8.20%  '

$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D35') '0.08194'
$ws.Range('E35').Value = '  +1.48%  '

$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D36') '0.03033'
$ws.Range('E36').Value = '  +4.28%  '

Set-TextValue $ws.Range('D37') '11.32'
$ws.Range('E37').Value = '  +8.99%  '

Set-TextValue $ws.Range('D38') '6.737'
$ws.Range('E38').Value = '  +0.40%  '

Set-TextValue $ws.Range('D39') '0.2731'
$ws.Range('E39').Value = '  +1.41%  '

Set-TextValue $ws.Range('D40') '0.09186'
$ws.Range('E40').Value = '  -0.63%  '

$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D41') '13.77'
$ws.Range('E41').Value = '  +5.07%  '

$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D42') '0.7638'
$ws.Range('E42').Value = '  +1.47%  '

Set-TextValue $ws.Range('D43') '1.428'
$ws.Range('E43').Value = '  -0.11%  '

Set-TextValue $ws.Range('D44') '16.60'
$ws.Range('E44').Value = '  +2.76%  '

Set-TextValue $ws.Range('D45') '0.7060'
$ws.Range('E45').Value = '  +2.07%  '

Set-TextValue $ws.Range('D46') '2.516'
$ws.Range('E46').Value = '  +2.24%  '

Set-TextValue $ws.Range('D48') '0.9999'
$ws.Range('E48').Value = '  -0.30%  '

Set-TextValue $ws.Range('D49') '0.08345'
$ws.Range('E49').Value = '  +0.07%  '

Set-TextValue $ws.Range('D50') '135.78'
$ws.Range('E50').Value = '  +1.89%  '

Set-TextValue $ws.Range('D51') '1.246'
$ws.Range('E51').Value = '  -1.60%  '
